$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.964.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.218.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.92%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.27"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.00"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -7.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.101"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.555.47"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.841"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.223.78"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.883.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0964"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.88"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.63"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.68"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.97"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.90%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +13.08%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.55"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0807"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.118"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.82"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.35"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.08"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.17%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -9.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.70"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.88"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.96%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.37%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -16.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.27%  "
